$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) from H1 into the two new header cells
# so I1/J1 get the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header labels for the new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-17
$data = @(
    @(7, 7),
    @(8, 8),
    @(11, 12),
    @(3, 4),
    @(7, 7),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(7, 7),
    @(7, 7),
    @(5, 6),
    @(7, 7),
    @(4, 4),
    @(4, 4),
    @(4, 4),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
